$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.724.35'
$ws.Range('E2').Value = '  -0.19%  '
$ws.Range('D3').Value = '3.150.31'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.56'
$ws.Range('E5').Value = '  +0.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.48'
$ws.Range('E6').Value = '  -1.84%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '3.149.42'
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('E9').Value = '  -0.50%  '
$ws.Range('E10').Value = '  -2.09%  '
$ws.Range('E11').Value = '  -1.08%  '
$ws.Range('E12').Value = '  -0.83%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000260'
$ws.Range('E13').Value = '  +2.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.04'
$ws.Range('E14').Value = '  -1.28%  '
$ws.Range('D15').Value = '3.668.46'
$ws.Range('E15').Value = '  +0.17%  '
$ws.Range('D16').Value = '64.767.90'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').Value = '3.150.58'
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('E18').Value = '  -1.65%  '
$ws.Range('E19').Value = '  +0.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '503.78'
$ws.Range('E20').Value = '  -1.57%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.83'
$ws.Range('E21').Value = '  -0.70%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.712'
$ws.Range('E23').Value = '  -3.28%  '
$ws.Range('E24').Value = '  -1.48%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.89'
$ws.Range('E25').Value = '  -1.36%  '
$ws.Range('E26').Value = '  -0.38%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.87'
$ws.Range('E27').Value = '  +1.86%  '
$ws.Range('E28').Value = '  -1.18%  '
$ws.Range('E29').Value = '  -1.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.80'
$ws.Range('E30').Value = '  +5.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '27.43'
$ws.Range('E31').Value = '  -1.93%  '
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.20'
$ws.Range('E33').Value = '  +1.04%  '
$ws.Range('E34').Value = '  +1.28%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.44'
$ws.Range('E35').Value = '  -1.94%  '
$ws.Range('E36').Value = '  -1.61%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0892'
$ws.Range('E37').Value = '  +3.42%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '475.47'
$ws.Range('E38').Value = '  -1.68%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0414'
$ws.Range('E39').Value = '  -2.34%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.92'
$ws.Range('E40').Value = '  -3.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.66'
$ws.Range('E41').Value = '  +0.31%  '
$ws.Range('D42').Value = '3.000.25'
$ws.Range('E42').Value = '  -3.77%  '
$ws.Range('E43').Value = '  -3.74%  '
$ws.Range('E44').Value = '  -2.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.280'
$ws.Range('E45').Value = '  -3.85%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '28.03'
$ws.Range('E46').Value = '  -3.92%  '
$ws.Range('D47').Value = '0.0₃0579'
$ws.Range('E47').Value = '  +1.03%  '
$ws.Range('E48').Value = '  -0.02%  '
$ws.Range('E49').Value = '  -1.62%  '
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.22'
$ws.Range('E50').Value = '  -2.87%  '
$ws.Range('B51').Value = 'Arweave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '33.81'
$ws.Range('E51').Value = '  +7.85%  '
